$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) and Volume(1h) (E) columns with the latest scrape.
# Numeric-looking new values (e.g. "111.96") are written with a leading
# apostrophe so Excel stores them as literal text -- matching the sheet's
# existing inline-string cells -- instead of silently converting them to
# numbers; the cell style is then reset to "Normal" so no stray quote-prefix
# / number-format style is left behind.

$ws.Range("D2").Value = '42.414.64'
$ws.Range("E2").Value = '  -2.94%  '
$ws.Range("D3").Value = '2.220.31'
$ws.Range("E3").Value = '  -2.60%  '
$ws.Range("E4").Value = '  +0.56%  '
$ws.Range("D5").Value = '''111.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -9.12%  '
$ws.Range("D6").Value = '''297.69'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +11.98%  '
$ws.Range("D7").Value = '''0.629'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.46%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Value = '''0.611'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.13%  '
$ws.Range("D10").Value = '''45.21'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.66%  '
$ws.Range("D11").Value = '''0.0927'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.15%  '
$ws.Range("D12").Value = '''54.87'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.02%  '
$ws.Range("D13").Value = '''8.88'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.01%  '
$ws.Range("E14").Value = '  -2.37%  '
$ws.Range("D15").Value = '''0.951'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.60%  '
$ws.Range("D16").Value = '''15.13'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.02%  '
$ws.Range("D17").Value = '2.558.92'
$ws.Range("E17").Value = '  -2.27%  '
$ws.Range("D18").Value = '2.227.78'
$ws.Range("E18").Value = '  -2.39%  '
$ws.Range("D19").Value = '42.363.36'
$ws.Range("E19").Value = '  -2.94%  '
$ws.Range("D20").Value = '''7.33'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.80%  '
$ws.Range("D21").Value = '''0.0000106'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.25%  '
$ws.Range("D22").Value = '''73.87'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.17%  '
$ws.Range("D23").Value = '''3.49'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +21.43%  '
$ws.Range("E24").Value = '  -5.84%  '
$ws.Range("D25").Value = '''229.63'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.76%  '
$ws.Range("D26").Value = '''9.39'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.43%  '
$ws.Range("D27").Value = '''11.74'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.82%  '
$ws.Range("D28").Value = '''0.999'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.78%  '
$ws.Range("D29").Value = '''3.91'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.46%  '
$ws.Range("D30").Value = '''38.75'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -8.58%  '
$ws.Range("E31").Value = '  -1.82%  '
$ws.Range("D33").Value = '''174.55'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.46%  '
$ws.Range("D34").Value = '''21.14'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.85%  '
$ws.Range("D35").Value = '''0.0887'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.41%  '
$ws.Range("D36").Value = '''5.71'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.73%  '
$ws.Range("D37").Value = '''4.94'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.47%  '
$ws.Range("D38").Value = '''4.28'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.15%  '
$ws.Range("E39").Value = '  -2.14%  '
$ws.Range("D40").Value = '''0.0368'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.75%  '
$ws.Range("E41").Value = '  -4.23%  '
$ws.Range("D42").Value = '''2.50'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.34%  '
$ws.Range("D43").Value = '''0.237'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.45%  '
$ws.Range("D44").Value = '''70.90'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.85%  '
$ws.Range("D45").Value = '''13.06'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.41%  '
$ws.Range("E46").Value = '  +0.23%  '
$ws.Range("E47").Value = '  -3.41%  '
$ws.Range("E48").Value = '  -3.72%  '
$ws.Range("E49").Value = '  +4.58%  '
$ws.Range("D50").Value = '''104.68'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.60%  '
$ws.Range("D51").Value = '''8.53'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.09%  '
